# Generate Report for Handoff
# Updates the localization-status workbook to reflect that "b.md" is now
# ready for handoff: status flips from "False" to "Ready for handoff",
# new handoff xliff files + datetimes are recorded per-locale, and an
# error detail note (stale handback version) is attached. The Overview
# sheet is refreshed to match, and the "Error Detail" column is widened.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8901f55738162e0cd0fdac90cb219d22c0b60c75/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8f9d9ef0513a8f1221130fc82bfff73f7553fbb9/e2e/b.md."

# ---- Overview sheet ------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-16 08:37:21"

# ---- zh-cn sheet -----------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-16 08:37:15"
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Columns.Item(16).ColumnWidth = 39.2

# ---- de-de sheet -----------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-16 08:37:21"
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = 39.2
